$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Extra product" sheet -> rename to "Scrap", move to the end, and replace
#    its single product row with the new Gelatin/Pectin/Tapioca JAR-style data.
#    (Its old row -- G50 Centrum Men / Centrum MultiGummies Men / Centrum men --
#    is relocated onto "Sample Info" as a new row 10, see step 2.)
# ---------------------------------------------------------------------------
$scrap = $wb.Worksheets.Item("Extra product")
$scrap.Name = "Scrap"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scrap.Move($null, $lastSheet)

# Re-fetch by name: after Move() any previously-held reference resolves by
# the sheet's *old* positional index rather than following the sheet itself.
$scrap = $wb.Worksheets.Item("Scrap")
$scrap.Cells.Clear()

$scrapData = @(
    @("Gelatin","Gelatin","Factor",0,"Absent"),
    @("Gelatin","Gelatin","Factor",1,"Present"),
    @("Pectin","Pectin","Factor",0,"Absent"),
    @("Pectin","Pectin","Factor",1,"Present"),
    @("Tapioca","Tapioca","Factor",0,"Absent"),
    @("Tapioca","Tapioca","Factor",1,"Present")
)
$r = 1
foreach ($row in $scrapData) {
    $scrap.Cells.Item($r, 1).Value = $row[0]
    $scrap.Cells.Item($r, 2).Value = $row[1]
    $scrap.Cells.Item($r, 3).Value = $row[2]
    $scrap.Cells.Item($r, 4).Value = $row[3]
    $scrap.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. "Sample Info" -- append the old "Extra product" row as row 10.
# ---------------------------------------------------------------------------
$sampleInfo = $wb.Worksheets.Item("Sample Info")
$sampleInfo.Cells.Item(10, 1).Value = "G50 Centrum Men"
$sampleInfo.Cells.Item(10, 2).Value = "Centrum MultiGummies Men"
$sampleInfo.Cells.Item(10, 3).Value = "Centrum men"

# ---------------------------------------------------------------------------
# 3. "Panel Details" -- renumber clusters, add Min/Max columns (D,E) with a
#    #,##0.00 number format on the "A_..." measurement rows, and append a new
#    "First Bite Hardness" row.
#    (Column/text writes below are ordered to match the shared-string table
#    insertion order of the target workbook: the new row 37 text comes before
#    the "Min"/"Max" headers.)
# ---------------------------------------------------------------------------
$panelDetails = $wb.Worksheets.Item("Panel Details")

$bigNum = [double]"9.9999999999999997E+98"
$negBigNum = [double]"-9.9999999999999997E+98"

$clusters = @(1,1,1,1,2,2,3,2,2,3,2,4,5,5,5,5,5,5,5,5,4,5,4,6,6,6,7,7,7,7,7,7,7,7,7)
for ($i = 0; $i -lt $clusters.Count; $i++) {
    $row = $i + 2
    $panelDetails.Cells.Item($row, 3).Value = $clusters[$i]
    $panelDetails.Cells.Item($row, 4).Value = 0
    $panelDetails.Cells.Item($row, 5).Value = 100
}

# Rows 28-36 (the "A_" texture-analyzer measurements) get real Min/Max bounds
# instead of the generic 0/100 placeholder, formatted as #,##0.00.
$panelDetails.Range("D28:E36").NumberFormat = "#,##0.00"

$panelDetails.Range("D28").Value = $negBigNum
$panelDetails.Range("E28").Value = $bigNum

$panelDetails.Range("D29").Value = 0
$panelDetails.Range("E29").Value = $bigNum

$panelDetails.Range("D30").Value = 0
$panelDetails.Range("E30").Value = $bigNum

$panelDetails.Range("D31").Value = 0
$panelDetails.Range("E31").Value = $bigNum

$panelDetails.Range("D32").Value = 0
$panelDetails.Range("E32").Value = $bigNum

$panelDetails.Range("D33").Value = $negBigNum
$panelDetails.Range("E33").Value = 0

$panelDetails.Range("D34").Value = $negBigNum
$panelDetails.Range("E34").Value = 0

$panelDetails.Range("D35").Value = 0
$panelDetails.Range("E35").Value = $bigNum

$panelDetails.Range("D36").Value = 0
$panelDetails.Range("E36").Value = $bigNum

# New row 37: "First Bite Hardness"  (new shared strings 181/182)
$panelDetails.Cells.Item(37, 1).Value = "First Bite Hardness"
$panelDetails.Cells.Item(37, 2).Value = "A_First Bite Hardness"
$panelDetails.Cells.Item(37, 3).Value = 7
$panelDetails.Range("D37").NumberFormat = "#,##0.00"
$panelDetails.Range("E37").NumberFormat = "#,##0.00"
$panelDetails.Range("D37").Value = 0
$panelDetails.Range("E37").Value = $bigNum

# New header cells D1/E1  (new shared strings 183/184 -- written last so they
# land after the Scrap-sheet and row-37 strings in the shared string table)
$panelDetails.Range("D1").Value = "Min"
$panelDetails.Range("E1").Value = "Max"

# Cosmetic column widths for the newly-touched columns.
$panelDetails.Columns.Item(4).ColumnWidth = 28.5
$panelDetails.Columns.Item(5).ColumnWidth = 43.67
$panelDetails.Columns.Item(6).ColumnWidth = 29
$panelDetails.Columns.Item(7).ColumnWidth = 23.17
$panelDetails.Columns.Item(8).ColumnWidth = 47
$panelDetails.Columns.Item(9).ColumnWidth = 29.67

$panelDetails.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Sheet view / selection bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$sampleInfo.Activate()
$sampleInfo.Range("B1:B10").Select() | Out-Null

$respondentDetails = $wb.Worksheets.Item("Respondent Details")
$respondentDetails.Activate()
$respondentDetails.Range("A69:XFD74").Select() | Out-Null

$scrap.Activate()
$scrap.Range("O38").Select() | Out-Null

$panelDetails.Activate()
$panelDetails.Range("F35").Select() | Out-Null
